$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = '56611'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2808'
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = '67650'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2528'
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = '50340'
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '45098'
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '10580'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '5709'
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '13468'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '5469'
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '17366'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '5201'
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '18596'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5129'
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '29561'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '4572'
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '53354'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2949'
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '67477'
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '13929'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '5432'
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '15948'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '5290'
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '16465'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '54698813'
$ws.Range("C18").Value = '閃亮唐老鴨'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '5256'
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = '16561'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '31495601'
$ws.Range("C19").Value = '陈晓军'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '5252'
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = '20437'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '54085771'
$ws.Range("C20").Value = '㊥Matthieu'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5029'
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = '21281'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '55769051'
$ws.Range("C21").Value = '㊥叮叮当.'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4985'
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = '23067'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '46289694'
$ws.Range("C22").Value = '㊥Vincent'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4884'
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = '30083'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4550'
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = '31175'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '56732705'
$ws.Range("C24").Value = '时间温柔皆遗憾'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '4505'
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = '31407'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '58839983'
$ws.Range("C25").Value = '每逢佳节胖六斤'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '4496'
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = '33215'
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '56585361'
$ws.Range("C26").Value = '"㊥ go策划我要ali"'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4414'
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = '39503'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '4108'
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = '39770'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '4097'
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = '44212'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '3729'
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = '5930'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '6139'
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = '8090'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '5917'
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = '11580'
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = '55317038'
$ws.Range("C32").Value = 'necman12345'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '5631'
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = '12071'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '5590'
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = '12344'
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = '11582001'
$ws.Range("C34").Value = 'iMinatoX4'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '5564'
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = '14313'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '5400'
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = '17893'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '5169'
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = '20033'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '5052'
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = '29670'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '4567'
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = '30904'
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = '56573048'
$ws.Range("C39").Value = 'Xiaotian'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4515'
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = '31223'
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = '47459684'
$ws.Range("C40").Value = '㊥阿闹切克闹'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4503'
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = '32669'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4439'
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = '33216'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4414'
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = '36512'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4256'
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = '39491'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4109'
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = '40573'
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = '58203298'
$ws.Range("C45").Value = '权旨qua'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '4051'
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = '41966'
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = '59020292'
$ws.Range("C46").Value = 'Sharnoth'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3990'
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = '42463'
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = '38893233'
$ws.Range("C47").Value = '"快乐 二哈"'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '3954'
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = '42825'
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = '32316256'
$ws.Range("C48").Value = '"秋の風 .."'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '3908'
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = '49131'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '3230'
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = '58033'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '2758'
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = '67347'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '2532'
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = '61638'
$ws.Range("E52").NumberFormat = "@"
$ws.Range("E52").Value = '2649'
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = '50508'
$ws.Range("E53").NumberFormat = "@"
$ws.Range("E53").Value = '3119'
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = '42412'
$ws.Range("E56").NumberFormat = "@"
$ws.Range("E56").Value = '3960'
$ws.Range("A57").NumberFormat = "@"
$ws.Range("A57").Value = '51032'
$ws.Range("E57").NumberFormat = "@"
$ws.Range("E57").Value = '3080'
$ws.Range("A58").NumberFormat = "@"
$ws.Range("A58").Value = '58739'
$ws.Range("E58").NumberFormat = "@"
$ws.Range("E58").Value = '2734'
$ws.Range("A59").NumberFormat = "@"
$ws.Range("A59").Value = '68005'
$ws.Range("E59").NumberFormat = "@"
$ws.Range("E59").Value = '2522'
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = '105082'
$ws.Range("E60").NumberFormat = "@"
$ws.Range("E60").Value = '1503'
$ws.Range("A61").NumberFormat = "@"
$ws.Range("A61").Value = '106741'
$ws.Range("E61").NumberFormat = "@"
$ws.Range("E61").Value = '1500'
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = '109637'
$ws.Range("E62").NumberFormat = "@"
$ws.Range("E62").Value = '1470'
$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = '121110'
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = '50398'
$ws.Range("E75").NumberFormat = "@"
$ws.Range("E75").Value = '3127'
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = '89788'
$ws.Range("E78").NumberFormat = "@"
$ws.Range("E78").Value = '1888'
$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = '96275'
$ws.Range("E79").NumberFormat = "@"
$ws.Range("E79").Value = '1635'
$ws.Range("A80").NumberFormat = "@"
$ws.Range("A80").Value = '158270'
$ws.Range("A81").NumberFormat = "@"
$ws.Range("A81").Value = '210360'
